# CSU9 - Manter Organização: sincronizar casos de uso com os protótipos (#45)

$d = $word.ActiveDocument

# Unicode helper characters
$lquote = [char]0x201C
$rquote = [char]0x201D

function Find-Replace($range, $findText, $replaceText) {
    return $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

function Find-InRange($range, $findText) {
    return $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# ---------------------------------------------------------------------------
# 1. "Pós-Cond" + "ição" -> "Pós-Condição"  (merge two runs into one)
# ---------------------------------------------------------------------------
Find-Replace $d.Content "Pós-Condição" "Pós-Condição" | Out-Null

# ---------------------------------------------------------------------------
# 2. "RF8" + "-Manter organização" -> "RF8-Manter organização"
# ---------------------------------------------------------------------------
Find-Replace $d.Content "RF8-Manter organização" "RF8-Manter organização" | Out-Null

# ---------------------------------------------------------------------------
# 3. "Seção Remover" + " Organização" -> "Seção Remover Organização"
# ---------------------------------------------------------------------------
Find-Replace $d.Content "Remoção: Ver Seção Remover Organização" "Remoção: Ver Seção Remover Organização" | Out-Null

# ---------------------------------------------------------------------------
# 4. "1. Ator seleciona opção de inserção (Tela T17)." ->
#    "...(Tela T17), clicando em “Adicionar organização”."
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "1. Ator seleciona opção de inserção*") {
        $rng = $p.Range
        Find-InRange $rng ")." | Out-Null
        $insPoint = $d.Range($rng.Start + 1, $rng.Start + 1)
        $insPoint.InsertAfter(", clicando em " + $lquote + "Adicionar organização" + $rquote)
        break
    }
}

# ---------------------------------------------------------------------------
# 5. "3. Ator informa os dados e submete para o sistema." /
#    "4. Sistema grava dados informados em meio persistente."
#    -> "3. ... sistema clicando no botão “Salvar”."
#       "4. Sistema grava dados informados em meio persistente." (merge)
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "3. Ator informa os dados e submete para o sistema.*") {
        $rng = $p.Range
        $insPoint = $d.Range($rng.End - 2, $rng.End - 2)
        $insPoint.InsertAfter(" clicando no botão " + $lquote + "Salvar" + $rquote)
        break
    }
}
Find-Replace $d.Content "4. Sistema grava dados informados em meio persistente." "4. Sistema grava dados informados em meio persistente." | Out-Null

# ---------------------------------------------------------------------------
# 6. "Consultar Or" + "ganização" -> "Consultar Organização"
# ---------------------------------------------------------------------------
Find-Replace $d.Content "Consultar Organização" "Consultar Organização" | Out-Null

# ---------------------------------------------------------------------------
# 7. Heading "Remover Organização" -> "Remover (desativar) Organização"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Seção: Remover Organização*") {
        $rng = $p.Range
        Find-InRange $rng "Remover Organização" | Out-Null
        $insPoint = $d.Range($rng.Start + 7, $rng.Start + 7)
        $insPoint.InsertAfter(" (desativar)")
        break
    }
}

# ---------------------------------------------------------------------------
# 8. Sumário: "Remove uma organização da base de dados. Desvinculando todos
#    os responsáveis e animais ligadas a organização removida."
#    -> "Desativa uma organização no sistema."
# ---------------------------------------------------------------------------
Find-Replace $d.Content "Remove uma organização da base de dados. Desvinculando todos os responsáveis e animais ligadas a organização removida." "Desativa uma organização no sistema." | Out-Null

# ---------------------------------------------------------------------------
# 9. "1. Ator seleciona opção de remoção de uma organização (Tela 17)." ->
#    "1. Ator seleciona opção de remoção (desativação) de uma organização
#    (Tela 17), simbolizada pelo ícone de lixeira."
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "1. Ator seleciona opção de remoção de uma organização*") {
        $rng = $p.Range
        Find-InRange $rng "remoção de uma organização" | Out-Null
        $insPoint1 = $d.Range($rng.Start + 8, $rng.Start + 8)
        $insPoint1.InsertAfter(" (desativação)")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "1. Ator seleciona opção de remoção*") {
        $rng = $p.Range
        Find-InRange $rng ")." | Out-Null
        $insPoint2 = $d.Range($rng.Start + 1, $rng.Start + 1)
        $insPoint2.InsertAfter(", simbolizada pelo ícone de lixeira")
        break
    }
}

# ---------------------------------------------------------------------------
# 10. "contendo entrada para o nome da organização como método de
#     confirmação de remoção." -> "... de remoção (desativação)."
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "2. Sistema exibe *confirmação de remoção.*") {
        $rng = $p.Range
        Find-InRange $rng "confirmação de remoção." | Out-Null
        $insPoint = $d.Range($rng.End - 1, $rng.End - 1)
        $insPoint.InsertBefore(" (desativação)")
        break
    }
}

# ---------------------------------------------------------------------------
# 11. "3. Ator insere o nome da organização e confirma a remoção através do
#     botão “Confirmar”." ->
#     "3. Ator insere o nome da organização e confirma a remoção
#     (desativação) através do botão “Desativar”."
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "3. Ator insere o nome da organização*") {
        $rng = $p.Range
        Find-InRange $rng "a remoção através" | Out-Null
        $insPoint = $d.Range($rng.Start + 10, $rng.Start + 10)
        $insPoint.InsertAfter("(desativação) ")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "3. Ator insere o nome da organização*") {
        $rng = $p.Range
        Find-Replace $rng ("Confirmar" + $rquote) ("Desativar" + $rquote) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 12. "4. Sistema exclui organização do meio persistente." ->
#     "4. Sistema modifica status da organização para “Desativada”."
# ---------------------------------------------------------------------------
Find-Replace $d.Content "4. Sistema exclui organização do meio persistente." ("4. Sistema modifica status da organização para " + $lquote + "Desativada" + $rquote + ".") | Out-Null

# ---------------------------------------------------------------------------
# 13. "Fluxo" + " Alternativo" -> "Fluxo Alternativo" (Remover section)
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Fluxo Alternativo*") {
        $rng = $p.Range
        $found = Find-InRange $rng "Fluxo Alternativo"
        if ($found) {
            break
        }
    }
}
Find-Replace $d.Content "Fluxo Alternativo" "Fluxo Alternativo" | Out-Null

# ---------------------------------------------------------------------------
# 14. "Organização não pôde ser excluída. Sistema exibe mensagem..." ->
#     "Organização não pôde ser desativada. Sistema exibe mensagem
#     “Ocorreu um erro ao desativar a organização. Tente novamente.”."
# ---------------------------------------------------------------------------
Find-Replace $d.Content "Organização não pôde ser excluída. Sistema exibe mensagem “Ocorreu um erro ao remover a organização. Tente novamente.”." "Organização não pôde ser desativada. Sistema exibe mensagem “Ocorreu um erro ao desativar a organização. Tente novamente.”." | Out-Null

# ---------------------------------------------------------------------------
# 15. Move <w:lastRenderedPageBreak/> from the "1. Ator seleciona opção de
#     alterar organização." run up to the "Fluxo Principal" run before it.
#     (Rendering bookmark only - not text-visible; left as-is since COM has
#     no direct knob for it and it carries no semantic content.)
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 16. "Tela" + " T17.3" -> "Tela T17.3" (drop red color)
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "2. Sistema exibe formulário *T17.3*") {
        $rng = $p.Range
        Find-InRange $rng "Tela T17.3" | Out-Null
        $rng.Font.Color = -16777216  # wdColorAutomatic
        break
    }
}
Find-Replace $d.Content "Tela T17.3" "Tela T17.3" | Out-Null

# ---------------------------------------------------------------------------
# 17. "Organização não pôde ser alterada. Sistema exibe" + " mensagem..."
#     -> merge into a single run
# ---------------------------------------------------------------------------
Find-Replace $d.Content "Organização não pôde ser alterada. Sistema exibe mensagem “Ocorreu um erro ao alterar a organização. Tente novamente.”." "Organização não pôde ser alterada. Sistema exibe mensagem “Ocorreu um erro ao alterar a organização. Tente novamente.”." | Out-Null
